$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 161 (this shifts the existing rows 161-242
# down to 162-243, matching the dimension growing from A1:R242 to A1:R243).
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row with the new weekly price-report entry.
$ws.Range("A161").Value = 5
$ws.Range("B161").Value = "Macroferia Regional de Talca"
$ws.Range("C161").Value = "Maule"
$ws.Range("D161").Value = 44572
$ws.Range("E161").Value = 7
$ws.Range("F161").Value = 100114014
$ws.Range("G161").Value = "Betarraga"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 4000
$ws.Range("K161").Value = 500
$ws.Range("L161").Value = 500
$ws.Range("M161").Value = 500
$ws.Range("N161").Value = "`$/paquete 5 unidades"
$ws.Range("O161").Value = "Región del Maule"
$ws.Range("P161").Value = 100
$ws.Range("Q161").Value = 5
$ws.Range("R161").Value = "Hortaliza"
